# Delete the "Library Management System" / "Sita Leela Manas Jagannath" row
# (worksheet row 25) from the "Form responses 1" table, shifting all
# subsequent rows up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(25).Delete()

# Reflect the user's post-edit view state (scroll / zoom / selection)
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("G8").Select()
